$wb = $excel.ActiveWorkbook

$wsTutoren = $wb.Worksheets.Item("Tutoren")

# Fill in the bookkeeping formulas for rows 13..21 (one row per tutor).
# D: # of PUE tutorials designed, E: # of PUE tutorials reviewed,
# F: # of HUE tutorials designed, G: # of HUE tutorials reviewed,
# H/I: carry-forward helper columns mirroring F/H respectively.
for ($r = 13; $r -le 21; $r++) {
    $wsTutoren.Cells.Item($r, 4).Formula  = "=COUNTIF(PUE!J`$3:J`$16,Tutoren!`$B$r)"
    $wsTutoren.Cells.Item($r, 5).Formula  = "=COUNTIF(PUE!K`$3:L`$16,`$B$r)"
    $wsTutoren.Cells.Item($r, 6).Formula  = "=COUNTIF(HUE!J`$3:J`$16,Tutoren!`$B$r)"

    $gRow = $r + 3
    $wsTutoren.Cells.Item($r, 7).Formula  = "=COUNTIF(HUE!K`$3:`$L$gRow,Tutoren!`$B$r)"

    $wsTutoren.Cells.Item($r, 8).Formula  = "=F$r"
    $wsTutoren.Cells.Item($r, 9).Formula  = "=H$r"
}

# Restore the navigation trail recorded in the saved view state: the user
# last looked at PUE!J3, then HUE!J13, and finally landed on Tutoren!D14,
# which is left as the active sheet/tab.
[void]$wb.Worksheets.Item("PUE").Range("J3").Select()
[void]$wb.Worksheets.Item("HUE").Range("J13").Select()
[void]$wsTutoren.Range("D14").Select()
